$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N2 must remain a text string (not get auto-converted to a date serial)
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 67378816.54000001
$ws.Range("P2").Value = 648841600.25
$ws.Range("Q2").Value = 590708417.79
$ws.Range("R2").Value = 80.6117342152
$ws.Range("S2").Value = 404321344.3
$ws.Range("T2").Value = 404321344.3
$ws.Range("U2").Value = 92.03401891590001
$ws.Range("V2").Value = 8145489.89
$ws.Range("W2").Value = 45910684.99
$ws.Range("X2").Value = -511491.1
$ws.Range("Y2").Value = 68698127.58
$ws.Range("Z2").Value = 68640217.65000001
$ws.Range("AA2").Value = 1261401.11
$ws.Range("AG2").Value = 479523.22
$ws.Range("AP2").Value = 96.645149183
$ws.Range("AQ2").Value = 2166.10202991144
$ws.Range("AR2").Value = 3705.76905415071
$ws.Range("AS2").Value = 54784800.87
$ws.Range("AT2").Value = 3549.21835901981
